$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 83333450
$ws.Range("I33").Value = 126.63636
$ws.Range("J33").Value = 1000000000
$ws.Range("K33").Value = 126.63636
$ws.Range("L33").Value = 1000000000
$ws.Range("M33").Value = 102.36364
$ws.Range("N33").Value = -1000000458
$ws.Range("H40").Value = 2685.1428
$ws.Range("I40").Value = 2632.6667
$ws.Range("K40").Value = 2632.6667
$ws.Range("M40").Value = -2457.6667
$ws.Range("H43").Value = 11595.8
$ws.Range("I43").Value = 12163
$ws.Range("K43").Value = 12163
$ws.Range("M43").Value = -12094
$ws.Range("H59").Value = 6249.5
$ws.Range("J59").Value = 6249.5
$ws.Range("L59").Value = 18748.5
$ws.Range("N59").Value = -19862.5
$ws.Range("H135").Value = 1494.0667
$ws.Range("I135").Value = 1494.0344
$ws.Range("K135").Value = 13446.3096
$ws.Range("M135").Value = -10911.3096
$ws.Range("H137").Value = 2375.6667
$ws.Range("I137").Value = 3029.0833
$ws.Range("J137").Value = 1504.4445
$ws.Range("K137").Value = 9087.249899999999
$ws.Range("L137").Value = 4513.333500000001
$ws.Range("M137").Value = -6537.249899999999
$ws.Range("N137").Value = -9613.333500000001
$ws.Range("H138").Value = 1803.0159
$ws.Range("I138").Value = 888.1
$ws.Range("J138").Value = 2228.558
$ws.Range("K138").Value = 2664.3
$ws.Range("L138").Value = 6685.674
$ws.Range("M138").Value = 2475.7
$ws.Range("N138").Value = -16965.674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3097.375
$ws.Range("I61").Value = 2731.5
$ws.Range("J61").Value = 4195
$ws.Range("K61").Value = 2731.5
$ws.Range("L61").Value = 4195
$ws.Range("M61").Value = -2519.5
$ws.Range("N61").Value = -4619
$ws.Range("H136").Value = 3097.375
$ws.Range("I136").Value = 2731.5
$ws.Range("J136").Value = 4195
$ws.Range("K136").Value = 8194.5
$ws.Range("L136").Value = 12585
$ws.Range("M136").Value = -5644.5
$ws.Range("N136").Value = -17685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2905.8408
$ws.Range("I134").Value = 2786.4722
$ws.Range("K134").Value = 8359.4166
$ws.Range("M134").Value = -5824.4166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1390.4286
$ws.Range("J31").Value = 981.3333
$ws.Range("L31").Value = 981.3333
$ws.Range("N31").Value = -1571.3333
$ws.Range("H34").Value = 1390.4286
$ws.Range("J34").Value = 981.3333
$ws.Range("L34").Value = 981.3333
$ws.Range("N34").Value = -1385.3333
$ws.Range("H62").Value = 10228.667
$ws.Range("J62").Value = 9368.125
$ws.Range("L62").Value = 9368.125
$ws.Range("N62").Value = -10616.125
$ws.Range("H65").Value = 10228.667
$ws.Range("J65").Value = 9368.125
$ws.Range("L65").Value = 46840.625
$ws.Range("N65").Value = -53080.625
$ws.Range("H132").Value = 2892.182
$ws.Range("I132").Value = 2626.25
$ws.Range("J132").Value = 3601.3333
$ws.Range("K132").Value = 7878.75
$ws.Range("L132").Value = 10803.9999
$ws.Range("M132").Value = -5348.75
$ws.Range("N132").Value = -15863.9999
$ws.Range("H134").Value = 2999.6667
$ws.Range("I134").Value = 2999
$ws.Range("K134").Value = 8997
$ws.Range("M134").Value = -6462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1118.3
$ws.Range("I5").Value = 1176
$ws.Range("J5").Value = 599
$ws.Range("K5").Value = 3528
$ws.Range("L5").Value = 1797
$ws.Range("M5").Value = -3416
$ws.Range("N5").Value = -2021
$ws.Range("H12").Value = 285.44446
$ws.Range("I12").Value = 312.81818
$ws.Range("K12").Value = 938.45454
$ws.Range("M12").Value = -765.45454
$ws.Range("H44").Value = 850
$ws.Range("I44").Value = 800
$ws.Range("K44").Value = 2400
$ws.Range("M44").Value = -2002
$ws.Range("H59").Value = 10749.5
$ws.Range("I59").Value = 11500
$ws.Range("J59").Value = 9999
$ws.Range("K59").Value = 34500
$ws.Range("L59").Value = 29997
$ws.Range("M59").Value = -33960
$ws.Range("N59").Value = -31077
$ws.Range("H113").Value = 2328
$ws.Range("I113").Value = 3096.25
$ws.Range("J113").Value = 1986.5555
$ws.Range("K113").Value = 9288.75
$ws.Range("L113").Value = 5959.666499999999
$ws.Range("M113").Value = -7118.75
$ws.Range("N113").Value = -10299.6665
$ws.Range("H135").Value = 1118.3
$ws.Range("I135").Value = 1176
$ws.Range("J135").Value = 599
$ws.Range("K135").Value = 10584
$ws.Range("L135").Value = 5391
$ws.Range("M135").Value = -8049
$ws.Range("N135").Value = -10461

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 204.41176
$ws.Range("J2").Value = 345.75
$ws.Range("L2").Value = 345.75
$ws.Range("N2").Value = -571.75
$ws.Range("H21").Value = 29499.666
$ws.Range("I21").Value = 62499
$ws.Range("J21").Value = 13000
$ws.Range("K21").Value = 62499
$ws.Range("L21").Value = 13000
$ws.Range("M21").Value = -62326
$ws.Range("N21").Value = -13346
$ws.Range("H22").Value = 2643.5557
$ws.Range("J22").Value = 799.3333
$ws.Range("L22").Value = 799.3333
$ws.Range("N22").Value = -1857.3333
$ws.Range("H30").Value = 29499.666
$ws.Range("I30").Value = 62499
$ws.Range("J30").Value = 13000
$ws.Range("K30").Value = 62499
$ws.Range("L30").Value = 13000
$ws.Range("M30").Value = -62394
$ws.Range("N30").Value = -13210
$ws.Range("H80").Value = 4200.769
$ws.Range("J80").Value = 2650.3333
$ws.Range("L80").Value = 2650.3333
$ws.Range("N80").Value = -4646.3333
$ws.Range("H83").Value = 4200.769
$ws.Range("J83").Value = 2650.3333
$ws.Range("L83").Value = 13251.6665
$ws.Range("N83").Value = -23235.6665
$ws.Range("H132").Value = 3100
$ws.Range("I132").Value = 2626.25
$ws.Range("J132").Value = 4995
$ws.Range("K132").Value = 7878.75
$ws.Range("L132").Value = 14985
$ws.Range("M132").Value = -5348.75
$ws.Range("N132").Value = -20045

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2537.7
$ws.Range("I22").Value = 940
$ws.Range("J22").Value = 2937.125
$ws.Range("K22").Value = 940
$ws.Range("L22").Value = 2937.125
$ws.Range("M22").Value = -645
$ws.Range("N22").Value = -3527.125
$ws.Range("H27").Value = 2537.7
$ws.Range("I27").Value = 940
$ws.Range("J27").Value = 2937.125
$ws.Range("K27").Value = 940
$ws.Range("L27").Value = 2937.125
$ws.Range("M27").Value = -833
$ws.Range("N27").Value = -3151.125
$ws.Range("H55").Value = 888.13336
$ws.Range("I55").Value = 452.6
$ws.Range("J55").Value = 1759.2
$ws.Range("K55").Value = 452.6
$ws.Range("L55").Value = 1759.2
$ws.Range("M55").Value = -279.6
$ws.Range("N55").Value = -2105.2
$ws.Range("H100").Value = 3752.2
$ws.Range("I100").Value = 3500
$ws.Range("J100").Value = 3860.2856
$ws.Range("K100").Value = 3500
$ws.Range("L100").Value = 3860.2856
$ws.Range("M100").Value = -2959
$ws.Range("N100").Value = -4942.2856
$ws.Range("H132").Value = 3353.125
$ws.Range("I132").Value = 2913.1428
$ws.Range("J132").Value = 6433
$ws.Range("K132").Value = 8739.428400000001
$ws.Range("L132").Value = 19299
$ws.Range("M132").Value = -6209.428400000001
$ws.Range("N132").Value = -24359
$ws.Range("H136").Value = 4926.6665
$ws.Range("I136").Value = 4926.6665
$ws.Range("K136").Value = 14779.9995
$ws.Range("M136").Value = -12229.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1972.3334
$ws.Range("I113").Value = 1544.2222
$ws.Range("J113").Value = 3256.6667
$ws.Range("K113").Value = 4632.6666
$ws.Range("L113").Value = 9770.000100000001
$ws.Range("M113").Value = -2462.6666
$ws.Range("N113").Value = -14110.0001
$ws.Range("H136").Value = 4745.484
$ws.Range("I136").Value = 2072.76
$ws.Range("J136").Value = 15881.833
$ws.Range("K136").Value = 6218.280000000001
$ws.Range("L136").Value = 47645.499
$ws.Range("M136").Value = -3668.280000000001
$ws.Range("N136").Value = -52745.499
